# Zerr_Suite.xlsx edit script
# Commit message: "Fixed a bug so that only rank 0 opens the output file.
#                   Also fixed the 1024 core zerr input"
#
# This applies two logical groups of changes to the "Zerr" worksheet:
#  1) The "1024 core zerr input" table (rows 19-21, 25, 28) is reworked:
#     row 19/20 switch from dividing by row12/row13 to dividing row6/row7
#     by 2*row21, and row21's column C becomes a hard-coded override (4
#     instead of the computed 8). Rows 25 and 28 are simplified to just
#     reference row 9 / row 10 respectively.
#  2) The sweep-time benchmark table (rows 42-49) gets corrected raw
#     timings for the 1024-core run (and the ones around it), a new label
#     cell "Cab" is added at D41, and the dependent speed-up formulas in
#     column E recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zerr")

# ---------------------------------------------------------------------
# Row 19 ("num_pin_x"): C19 becomes a literal override; D19:S19 now
# compute row6 / (2 * row21) instead of row6 / row12.
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 4
$ws.Range("D19:S19").Formula = "=D6/(2*D21)"

# ---------------------------------------------------------------------
# Row 20 ("num_pin_y"): C20 becomes a literal override; D20:S20 now
# compute row7 / (2 * row21) instead of row7 / row13.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20:S20").Formula = "=D7/(2*D21)"

# ---------------------------------------------------------------------
# Row 21 ("refinement"): C21 becomes a hard-coded 4 (was
# =MIN(C12:C13)/2, which evaluated to 8). D21:S21 keep their formula.
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 4

# ---------------------------------------------------------------------
# Row 25 ("num_cellsets_x"): simplified from row6/(2*row21) to just =row9
# ---------------------------------------------------------------------
$ws.Range("C25").Formula = "=C9"
$ws.Range("D25:S25").Formula = "=D9"

# ---------------------------------------------------------------------
# Row 28 ("num_cellsets_y"): simplified from row7/(2*row21) to just =row10
# ---------------------------------------------------------------------
$ws.Range("C28").Formula = "=C10"
$ws.Range("D28:S28").Formula = "=D10"

# ---------------------------------------------------------------------
# New label cell: D41 = "Cab" (adds a new shared string)
# ---------------------------------------------------------------------
$ws.Range("D41").Value = "Cab"

# ---------------------------------------------------------------------
# Sweep-time benchmark raw timings (column D, rows 42-49) - corrected
# numbers, most notably the fix for the 1024-core run (D49).
# Column E (speed-up ratios) recomputes automatically from these.
# ---------------------------------------------------------------------
$ws.Range("D42").Value = 0.30998799999999999
$ws.Range("D43").Value = 0.309363
$ws.Range("D44").Value = 0.309722
$ws.Range("D45").Value = 0.31451800000000002
$ws.Range("D46").Value = 0.33661000000000002
$ws.Range("D47").Value = 0.341335
$ws.Range("D48").Value = 0.36446200000000001
$ws.Range("D49").Value = 1.63008

# ---------------------------------------------------------------------
# View state: scroll back so row 1 is at the top again and select a
# single cell (X8) instead of the previous C28:S28 range selection.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("X8").Select()

$wb.Save()
